# Added the "add course" button to the GUI and it works.
# Record a first (real) student into the "students" sheet's first data row
# and move the selection to reflect where the user's cursor ended up
# after using the new Add-Course / Add-Student flow.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("students")

# --- Row 2: replace the placeholder "h" values with the first real
#     student record ---------------------------------------------------

# Plain text fields - Excel would not reinterpret these as numbers/dates,
# so a straight Value assignment is enough and keeps the existing style.
$ws.Range("B2").Value = "Hasith"
$ws.Range("C2").Value = "Dewmina"
$ws.Range("E2").Value = "M"
$ws.Range("F2").Value = "Lesli Kumara"
$ws.Range("H2").Value = "Ruwi, Muscat, Oman"

# Fields that look numeric/date-like ("1", "01/18/2006", "3456789") need
# to stay plain text (they are IDs / phone numbers / a typed date string,
# not real numbers), so enter them the same way a user would force text
# in Excel - with a leading apostrophe - rather than letting AutoDetect
# turn them into a number or a date serial.
$ws.Range("A2").Formula = "'1"
$ws.Range("D2").Formula = "'01/18/2006"
$ws.Range("G2").Formula = "'3456789"

# --- Move the active selection from H10 to H8 on the students sheet ---
[void]$ws.Activate()
[void]$ws.Range("H8").Select()
